$wb = $excel.ActiveWorkbook

# --- Sheet: private ---
$ws = $wb.Worksheets.Item("private")
$ws.Columns.Item(7).ColumnWidth = 10.8   # -> stored width ~11.711 (closest achievable)
$ws.Cells.Item(38, 7).Value = 264000.0
$ws.Cells.Item(38, 25).Value = -259256.6
$ws.Cells.Item(39, 7).Value = 290400.0
$ws.Cells.Item(39, 25).Value = -267563.8
$ws.Cells.Item(40, 7).Value = 319440.0
$ws.Cells.Item(40, 25).Value = -276701.5
$ws.Cells.Item(41, 7).Value = 351384.0
$ws.Cells.Item(41, 25).Value = -286753.3
$ws.Cells.Item(42, 7).Value = 386522.4
$ws.Cells.Item(42, 25).Value = -297810.5
$ws.Cells.Item(43, 7).Value = 425174.75
$ws.Cells.Item(43, 25).Value = -309973.75
$ws.Cells.Item(44, 7).Value = 467692.5
$ws.Cells.Item(44, 25).Value = -323351.4
$ws.Cells.Item(45, 7).Value = 514461.75
$ws.Cells.Item(45, 25).Value = -338067.75
$ws.Cells.Item(46, 7).Value = 565908.2
$ws.Cells.Item(46, 25).Value = -354255.9
$ws.Cells.Item(47, 7).Value = 622498.8
$ws.Cells.Item(47, 25).Value = -372748.4075
$ws.Cells.Item(48, 7).Value = 684748.9
$ws.Cells.Item(48, 25).Value = -394104.5675
$ws.Cells.Item(49, 7).Value = 753223.9
$ws.Cells.Item(49, 25).Value = -417595.575
$ws.Cells.Item(50, 7).Value = 828546.4
$ws.Cells.Item(50, 25).Value = -443437.1125
$ws.Cells.Item(51, 7).Value = 911401.15
$ws.Cells.Item(51, 25).Value = -471862.585
$ws.Cells.Item(52, 7).Value = 1002541.1
$ws.Cells.Item(52, 25).Value = -503130.4375
$ws.Cells.Item(53, 7).Value = 65285.5
$ws.Cells.Item(53, 25).Value = -1575033.9625
$ws.Cells.Item(54, 7).Value = 71814.5
$ws.Cells.Item(54, 25).Value = -1716619.0
$ws.Cells.Item(55, 7).Value = 78995.5
$ws.Cells.Item(55, 25).Value = -1872362.23
$ws.Cells.Item(56, 7).Value = 86895.5
$ws.Cells.Item(56, 25).Value = -2043679.9825
$ws.Cells.Item(57, 7).Value = 194585.0
$ws.Cells.Item(57, 25).Value = -2133129.02
$ws.Cells.Item(58, 7).Value = 214043.5
$ws.Cells.Item(58, 25).Value = -2155049.9625
$ws.Cells.Item(59, 7).Value = 235447.5
$ws.Cells.Item(59, 25).Value = -2373954.6475
$ws.Cells.Item(60, 7).Value = 258992.5
$ws.Cells.Item(60, 25).Value = -2614750.855
$ws.Cells.Item(61, 7).Value = 284891.4
$ws.Cells.Item(61, 25).Value = -2879625.8525
$ws.Cells.Item(62, 7).Value = 313380.6
$ws.Cells.Item(62, 25).Value = -3170987.83
$ws.Cells.Item(63, 7).Value = 344718.55
$ws.Cells.Item(63, 25).Value = -3491486.72
$ws.Cells.Item(64, 7).Value = 379190.45
$ws.Cells.Item(64, 25).Value = -3844035.57
$ws.Cells.Item(65, 7).Value = 417109.8
$ws.Cells.Item(65, 25).Value = -4231839.4825
$ws.Cells.Item(66, 7).Value = 458821.05
$ws.Cells.Item(66, 25).Value = -4658423.6
$ws.Cells.Item(67, 7).Value = 504702.7
$ws.Cells.Item(67, 25).Value = -5127667.295
$ws.Cells.Item(68, 7).Value = 324841.05
$ws.Cells.Item(68, 25).Value = -5874166.715
$ws.Cells.Item(69, 7).Value = 357325.35
$ws.Cells.Item(69, 25).Value = -6464982.6375
$ws.Cells.Item(70, 7).Value = 393057.95
$ws.Cells.Item(70, 25).Value = -7114881.165
$ws.Cells.Item(71, 7).Value = 432363.5
$ws.Cells.Item(71, 25).Value = -7829768.8675
$ws.Cells.Item(72, 7).Value = 475599.25
$ws.Cells.Item(72, 25).Value = -8616146.795

# --- Sheet: Income ---
$ws = $wb.Worksheets.Item("Income")
$ws.Columns.Item(7).ColumnWidth = 10.8   # -> stored width ~11.711 (closest achievable)
$ws.Cells.Item(6, 8).Value = 0.5
$ws.Cells.Item(7, 8).Value = 0.5
$ws.Cells.Item(8, 8).Value = 0.5
$ws.Cells.Item(9, 8).Value = 0.5
$ws.Cells.Item(10, 8).Value = 0.5
$ws.Cells.Item(11, 8).Value = 0.5
$ws.Cells.Item(12, 8).Value = 0.5
$ws.Cells.Item(13, 8).Value = 0.5
$ws.Cells.Item(14, 8).Value = 0.5
$ws.Cells.Item(15, 8).Value = 0.5
$ws.Cells.Item(16, 8).Value = 0.5
$ws.Cells.Item(17, 8).Value = 0.5
$ws.Cells.Item(18, 8).Value = 0.5
$ws.Cells.Item(19, 8).Value = 0.5
$ws.Cells.Item(20, 8).Value = 0.5
$ws.Cells.Item(21, 8).Value = 0.5
$ws.Cells.Item(22, 8).Value = 0.5
$ws.Cells.Item(23, 8).Value = 0.5
$ws.Cells.Item(24, 8).Value = 0.5
$ws.Cells.Item(25, 8).Value = 0.5
$ws.Cells.Item(26, 8).Value = 0.5
$ws.Cells.Item(27, 8).Value = 0.5
$ws.Cells.Item(28, 8).Value = 0.5
$ws.Cells.Item(29, 8).Value = 0.5
$ws.Cells.Item(30, 8).Value = 0.5
$ws.Cells.Item(31, 8).Value = 0.5
$ws.Cells.Item(32, 8).Value = 0.5
$ws.Cells.Item(33, 8).Value = 0.5
$ws.Cells.Item(34, 8).Value = 0.5
$ws.Cells.Item(35, 8).Value = 0.5
$ws.Cells.Item(36, 8).Value = 0.5
$ws.Cells.Item(37, 8).Value = 0.5
$ws.Cells.Item(38, 8).Value = 0.5
$ws.Cells.Item(39, 8).Value = 0.5
$ws.Cells.Item(40, 8).Value = 0.5
$ws.Cells.Item(41, 8).Value = 0.5
$ws.Cells.Item(42, 8).Value = 0.5
$ws.Cells.Item(43, 8).Value = 0.5
$ws.Cells.Item(44, 8).Value = 0.5
$ws.Cells.Item(45, 8).Value = 0.5
$ws.Cells.Item(46, 8).Value = 0.5
$ws.Cells.Item(47, 8).Value = 0.5
$ws.Cells.Item(48, 8).Value = 0.5
$ws.Cells.Item(49, 8).Value = 0.5
$ws.Cells.Item(50, 8).Value = 0.5
$ws.Cells.Item(51, 8).Value = 0.5
$ws.Cells.Item(52, 8).Value = 0.5
$ws.Cells.Item(53, 8).Value = 0.5
$ws.Cells.Item(54, 8).Value = 0.5
$ws.Cells.Item(55, 8).Value = 0.5
$ws.Cells.Item(56, 8).Value = 0.5
$ws.Cells.Item(57, 8).Value = 0.5
$ws.Cells.Item(58, 8).Value = 0.5
$ws.Cells.Item(59, 8).Value = 0.5
$ws.Cells.Item(60, 8).Value = 0.5
$ws.Cells.Item(61, 8).Value = 0.5
$ws.Cells.Item(62, 8).Value = 0.5
$ws.Cells.Item(63, 8).Value = 0.5
$ws.Cells.Item(64, 8).Value = 0.5
$ws.Cells.Item(65, 8).Value = 0.5
$ws.Cells.Item(66, 8).Value = 0.5
$ws.Cells.Item(67, 8).Value = 0.5
$ws.Cells.Item(68, 8).Value = 0.5
$ws.Cells.Item(69, 8).Value = 0.5
$ws.Cells.Item(70, 8).Value = 0.5
$ws.Cells.Item(71, 8).Value = 0.5
$ws.Cells.Item(72, 8).Value = 0.5
$ws.Cells.Item(38, 7).Value = 264000.0
$ws.Cells.Item(38, 25).Value = 66000.0
$ws.Cells.Item(39, 7).Value = 290400.0
$ws.Cells.Item(39, 25).Value = 72600.0
$ws.Cells.Item(40, 7).Value = 319440.0
$ws.Cells.Item(40, 25).Value = 79860.0
$ws.Cells.Item(41, 7).Value = 351384.0
$ws.Cells.Item(41, 25).Value = 87846.0
$ws.Cells.Item(42, 7).Value = 386522.4
$ws.Cells.Item(42, 25).Value = 96630.6
$ws.Cells.Item(43, 7).Value = 425174.75
$ws.Cells.Item(43, 25).Value = 106293.55
$ws.Cells.Item(44, 7).Value = 467692.5
$ws.Cells.Item(44, 25).Value = 116923.4
$ws.Cells.Item(45, 7).Value = 514461.75
$ws.Cells.Item(45, 25).Value = 128615.85
$ws.Cells.Item(46, 7).Value = 565908.2
$ws.Cells.Item(46, 25).Value = 141477.6
$ws.Cells.Item(47, 7).Value = 622498.8
$ws.Cells.Item(47, 25).Value = 155624.7
$ws.Cells.Item(48, 7).Value = 684748.9
$ws.Cells.Item(48, 25).Value = 171187.5
$ws.Cells.Item(49, 7).Value = 753223.9
$ws.Cells.Item(49, 25).Value = 188306.8
$ws.Cells.Item(50, 7).Value = 828546.4
$ws.Cells.Item(50, 25).Value = 207137.7
$ws.Cells.Item(51, 7).Value = 911401.15
$ws.Cells.Item(51, 25).Value = 227851.25
$ws.Cells.Item(52, 7).Value = 1002541.1
$ws.Cells.Item(52, 25).Value = 250636.1
$ws.Cells.Item(53, 7).Value = 65285.5
$ws.Cells.Item(53, 25).Value = -761810.0
$ws.Cells.Item(54, 7).Value = 71814.5
$ws.Cells.Item(54, 25).Value = -837991.1
$ws.Cells.Item(55, 7).Value = 78995.5
$ws.Cells.Item(55, 25).Value = -921791.1
$ws.Cells.Item(56, 7).Value = 86895.5
$ws.Cells.Item(56, 25).Value = -1013970.2
$ws.Cells.Item(57, 7).Value = 95585.0
$ws.Cells.Item(57, 25).Value = -1115367.6
$ws.Cells.Item(58, 7).Value = 105143.5
$ws.Cells.Item(58, 25).Value = -1226904.8
$ws.Cells.Item(59, 7).Value = 115657.5
$ws.Cells.Item(59, 25).Value = -1349595.3
$ws.Cells.Item(60, 7).Value = 127223.5
$ws.Cells.Item(60, 25).Value = -1484554.8
$ws.Cells.Item(61, 7).Value = 139945.5
$ws.Cells.Item(61, 25).Value = -1633010.3
$ws.Cells.Item(62, 7).Value = 153940.0
$ws.Cells.Item(62, 25).Value = -1796311.6
$ws.Cells.Item(63, 7).Value = 169334.0
$ws.Cells.Item(63, 25).Value = -1975943.2
$ws.Cells.Item(64, 7).Value = 186267.5
$ws.Cells.Item(64, 25).Value = -2173537.2
$ws.Cells.Item(65, 7).Value = 204894.5
$ws.Cells.Item(65, 25).Value = -2390891.0
$ws.Cells.Item(66, 7).Value = 225384.0
$ws.Cells.Item(66, 25).Value = -2629980.6
$ws.Cells.Item(67, 7).Value = 247922.0
$ws.Cells.Item(67, 25).Value = -2892979.5
$ws.Cells.Item(68, 7).Value = 42382.5
$ws.Cells.Item(68, 25).Value = -3412609.7
$ws.Cells.Item(69, 7).Value = 46621.0
$ws.Cells.Item(69, 25).Value = -3753870.2
$ws.Cells.Item(70, 7).Value = 51283.0
$ws.Cells.Item(70, 25).Value = -4129257.1
$ws.Cells.Item(71, 7).Value = 56411.0
$ws.Cells.Item(71, 25).Value = -4542183.0
$ws.Cells.Item(72, 7).Value = 62051.5
$ws.Cells.Item(72, 25).Value = -4996401.9

# --- Sheet: Folketrygden ---
$ws = $wb.Worksheets.Item("Folketrygden")
$ws.Columns.Item(7).ColumnWidth = 8.5   # -> stored width ~9.283 (closest achievable)
$ws.Cells.Item(6, 8).Value = 0.5
$ws.Cells.Item(7, 8).Value = 0.5
$ws.Cells.Item(8, 8).Value = 0.5
$ws.Cells.Item(9, 8).Value = 0.5
$ws.Cells.Item(10, 8).Value = 0.5
$ws.Cells.Item(11, 8).Value = 0.5
$ws.Cells.Item(12, 8).Value = 0.5
$ws.Cells.Item(13, 8).Value = 0.5
$ws.Cells.Item(14, 8).Value = 0.5
$ws.Cells.Item(15, 8).Value = 0.5
$ws.Cells.Item(16, 8).Value = 0.5
$ws.Cells.Item(17, 8).Value = 0.5
$ws.Cells.Item(18, 8).Value = 0.5
$ws.Cells.Item(19, 8).Value = 0.5
$ws.Cells.Item(20, 8).Value = 0.5
$ws.Cells.Item(21, 8).Value = 0.5
$ws.Cells.Item(22, 8).Value = 0.5
$ws.Cells.Item(23, 8).Value = 0.5
$ws.Cells.Item(24, 8).Value = 0.5
$ws.Cells.Item(25, 8).Value = 0.5
$ws.Cells.Item(26, 8).Value = 0.5
$ws.Cells.Item(27, 8).Value = 0.5
$ws.Cells.Item(28, 8).Value = 0.5
$ws.Cells.Item(29, 8).Value = 0.5
$ws.Cells.Item(30, 8).Value = 0.5
$ws.Cells.Item(31, 8).Value = 0.5
$ws.Cells.Item(32, 8).Value = 0.5
$ws.Cells.Item(33, 8).Value = 0.5
$ws.Cells.Item(34, 8).Value = 0.5
$ws.Cells.Item(35, 8).Value = 0.5
$ws.Cells.Item(36, 8).Value = 0.5
$ws.Cells.Item(37, 8).Value = 0.5
$ws.Cells.Item(38, 8).Value = 0.5
$ws.Cells.Item(39, 8).Value = 0.5
$ws.Cells.Item(40, 8).Value = 0.5
$ws.Cells.Item(41, 8).Value = 0.5
$ws.Cells.Item(42, 8).Value = 0.5
$ws.Cells.Item(43, 8).Value = 0.5
$ws.Cells.Item(44, 8).Value = 0.5
$ws.Cells.Item(45, 8).Value = 0.5
$ws.Cells.Item(46, 8).Value = 0.5
$ws.Cells.Item(47, 8).Value = 0.5
$ws.Cells.Item(48, 8).Value = 0.5
$ws.Cells.Item(49, 8).Value = 0.5
$ws.Cells.Item(50, 8).Value = 0.5
$ws.Cells.Item(51, 8).Value = 0.5
$ws.Cells.Item(52, 8).Value = 0.5
$ws.Cells.Item(53, 8).Value = 0.5
$ws.Cells.Item(54, 8).Value = 0.5
$ws.Cells.Item(55, 8).Value = 0.5
$ws.Cells.Item(56, 8).Value = 0.5
$ws.Cells.Item(57, 8).Value = 0.5
$ws.Cells.Item(58, 8).Value = 0.5
$ws.Cells.Item(59, 8).Value = 0.5
$ws.Cells.Item(60, 8).Value = 0.5
$ws.Cells.Item(61, 8).Value = 0.5
$ws.Cells.Item(62, 8).Value = 0.5
$ws.Cells.Item(63, 8).Value = 0.5
$ws.Cells.Item(64, 8).Value = 0.5
$ws.Cells.Item(65, 8).Value = 0.5
$ws.Cells.Item(66, 8).Value = 0.5
$ws.Cells.Item(67, 8).Value = 0.5
$ws.Cells.Item(68, 8).Value = 0.5
$ws.Cells.Item(69, 8).Value = 0.5
$ws.Cells.Item(70, 8).Value = 0.5
$ws.Cells.Item(71, 8).Value = 0.5
$ws.Cells.Item(72, 8).Value = 0.5
$ws.Cells.Item(57, 7).Value = 99000.0
$ws.Cells.Item(57, 25).Value = 99000.0
$ws.Cells.Item(58, 7).Value = 108900.0
$ws.Cells.Item(58, 25).Value = 108900.0
$ws.Cells.Item(59, 7).Value = 119790.0
$ws.Cells.Item(59, 25).Value = 119790.0
$ws.Cells.Item(60, 7).Value = 131769.0
$ws.Cells.Item(60, 25).Value = 131769.0
$ws.Cells.Item(61, 7).Value = 144945.9
$ws.Cells.Item(61, 25).Value = 144945.9
$ws.Cells.Item(62, 7).Value = 159440.6
$ws.Cells.Item(62, 25).Value = 159440.6
$ws.Cells.Item(63, 7).Value = 175384.55
$ws.Cells.Item(63, 25).Value = 175384.55
$ws.Cells.Item(64, 7).Value = 192922.95
$ws.Cells.Item(64, 25).Value = 192922.95
$ws.Cells.Item(65, 7).Value = 212215.3
$ws.Cells.Item(65, 25).Value = 212215.3
$ws.Cells.Item(66, 7).Value = 233437.05
$ws.Cells.Item(66, 25).Value = 233437.05
$ws.Cells.Item(67, 7).Value = 256780.7
$ws.Cells.Item(67, 25).Value = 256780.7
$ws.Cells.Item(68, 7).Value = 282458.55
$ws.Cells.Item(68, 25).Value = 282458.55
$ws.Cells.Item(69, 7).Value = 310704.35
$ws.Cells.Item(69, 25).Value = 310704.35
$ws.Cells.Item(70, 7).Value = 341774.95
$ws.Cells.Item(70, 25).Value = 341774.95
$ws.Cells.Item(71, 7).Value = 375952.5
$ws.Cells.Item(71, 25).Value = 375952.5
$ws.Cells.Item(72, 7).Value = 413547.75
$ws.Cells.Item(72, 25).Value = 413547.75
